# Replace the old 4-column roles/areas table with the new 14-column
# employee/inventory roster exported from the source system.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("A1").Value = "Número"
$ws.Range("B1").Value = "Nombre"
$ws.Range("C1").Value = "Activity"
$ws.Range("D1").Value = "categoria"
$ws.Range("E1").Value = "Confidencialidad"
$ws.Range("F1").Value = "Puesto"
$ws.Range("G1").Value = "Cód Centro Costo"
$ws.Range("H1").Value = "Departamento"
$ws.Range("I1").Value = "Supervisor"
$ws.Range("J1").Value = "Location"
$ws.Range("K1").Value = "Oracle Location"
$ws.Range("L1").Value = "Shift"
$ws.Range("M1").Value = "Responsabilidad Oracle"
$ws.Range("N1").Value = "Comentarios / Vacaciones"

# Data row 2
$ws.Range("A2").Value = 20687
$ws.Range("B2").Value = "ABREGO MORIN, CESAR EDUARDO"
$ws.Range("C2").Value = "Counting"
$ws.Range("D2").Value = "HOURLY DIRECT"
$ws.Range("E2").Value = "DIRECTO"
$ws.Range("F2").Value = "OPERADOR DE ENSAMBLE"
$ws.Range("G2").Value = 5001
$ws.Range("H2").Value = "Manufacturing Fabrication -Fab"
$ws.Range("I2").Value = "JESUS SERVANDO MUNGUIA MONTELO"
$ws.Range("J2").Value = "DH L2"
$ws.Range("K2").Value = "MONDP"
$ws.Range("L2").Value = 3

# Data row 3
$ws.Range("A3").Value = 19976
$ws.Range("B3").Value = "ABUNDIS GUAJARDO, REBECA"
$ws.Range("C3").Value = "Inventory Audit"
$ws.Range("D3").Value = "HOURLY INDIRECT"
$ws.Range("E3").Value = "DIRECTO"
$ws.Range("F3").Value = "INSPECTOR DE CALIDAD"
$ws.Range("G3").Value = 5226
$ws.Range("H3").Value = "Mfg OH -Quality / Continuous I"
$ws.Range("I3").Value = "JUAN EDGAR GUADALUPE GARCIA RU"
$ws.Range("J3").Value = "SKID"
$ws.Range("K3").Value = "MONSKID"
$ws.Range("L3").Value = 0

# Data row 4
$ws.Range("A4").Value = 21142
$ws.Range("B4").Value = "ABUNDIS VILLASANA, OMAR"
$ws.Range("C4").Value = "Counting"
$ws.Range("D4").Value = "HOURLY DIRECT"
$ws.Range("E4").Value = "DIRECTO"
$ws.Range("F4").Value = "OPERADOR DE ENSAMBLE"
$ws.Range("G4").Value = 5001
$ws.Range("H4").Value = "Manufacturing Fabrication -Fab"
$ws.Range("I4").Value = "FRANCISCO  ALEJANDRO MONTOYA C"
$ws.Range("J4").Value = "DH L2"
$ws.Range("K4").Value = "MONDP"
$ws.Range("L4").Value = 2

# Data row 5
$ws.Range("A5").Value = 21955
$ws.Range("B5").Value = "ABURTO BANDALA, VICTOR MANUEL"
$ws.Range("C5").Value = "Inventory Audit"
$ws.Range("D5").Value = "SALARY EXEMPT"
$ws.Range("E5").Value = "SALARY"
$ws.Range("F5").Value = "OPS PM PROJECT MANAGER"
$ws.Range("G5").Value = 6640
$ws.Range("H5").Value = "Global Operations"
$ws.Range("I5").Value = "KURC MARIAN"
$ws.Range("J5").Value = "MC"
$ws.Range("K5").Value = "MONMC"
$ws.Range("L5").Value = 0

# Autofit every populated column to its contents (A through N)
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(8).AutoFit()
$ws.Columns.Item(9).AutoFit()
$ws.Columns.Item(10).AutoFit()
$ws.Columns.Item(11).AutoFit()
$ws.Columns.Item(12).AutoFit()
$ws.Columns.Item(13).AutoFit()
$ws.Columns.Item(14).AutoFit()

# Leave the same cell selected as in the authored workbook
[void]$ws.Range("J1").Select()
